$wb = $excel.ActiveWorkbook

# "constant_vol_surface" is the first sheet (sheetId=2, tab-selected) and
# holds the flat 10% vol curve used by the new time-dependent GBM pricer
# test. Bump the vols from 10 -> 40 across the whole B2:B13 column.
$ws = $wb.Worksheets.Item("constant_vol_surface")
$ws.Activate()
$ws.Range("B2:B13").Value = 40

# Leave the selection where the author left it when they saved the file.
$ws.Range("I9").Select()
